# Update countries & provincias Spain
# - Refresh the "last updated" timestamp.
# - Refresh case numbers for India, Pakistan, Tailandia, Islas Turcas y
#   Caicos and Mongolia.
# - Montserrat's case count overtakes Islas Malvinas, so the two rows swap
#   places (country name + stats) while every other row stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados a ..." timestamp ---
$ws.Range("A1").Value() = "Datos actualizados a 18 de Septiembre de 2020 a las 06:22"

# --- India (row 5) ---
$ws.Range("B5").Value() = 5214677
$ws.Range("C5").Value() = 1991
$ws.Range("D5").Value() = 4112551
$ws.Range("E5").Value() = 1017722

# --- Pakistan (row 21) ---
$ws.Range("B21").Value() = 304386
$ws.Range("C21").Value() = 752
$ws.Range("D21").Value() = 291683
$ws.Range("E21").Value() = 6295
$ws.Range("G21").Value() = 9
$ws.Range("H21").Value() = 6408

# --- Tailandia (row 133) ---
$ws.Range("B133").Value() = 3497
$ws.Range("C133").Value() = 7
$ws.Range("D133").Value() = 3328
$ws.Range("E133").Value() = 111

# --- Islas Turcas y Caicos (row 172) ---
$ws.Range("B172").Value() = 663
$ws.Range("C172").Value() = 4
$ws.Range("D172").Value() = 565

# --- Mongolia (row 185) ---
$ws.Range("D185").Value() = 302
$ws.Range("E185").Value() = 9

# --- Montserrat / Islas Malvinas (rows 214-215) swap places ---
# Row 214 becomes Montserrat's (updated) data, row 215 becomes what used
# to be Montserrat's neighbour, Islas Malvinas, with its original data.
$ws.Range("A214").Value() = "Montserrat"
$ws.Range("D214").Value() = 12
$ws.Range("H214").Value() = 1

$ws.Range("A215").Value() = "Islas Malvinas"
$ws.Range("D215").Value() = 13
$ws.Range("H215").Value() = 0
